$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Leaderboard" -> strike-through (paragraph mark + run)
# ---------------------------------------------------------------------
$leaderPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Leaderboard*") {
        $leaderPara = $p
        break
    }
}
$leaderPara.Range.Font.StrikeThrough = $true

# ---------------------------------------------------------------------
# 2) "Camera should be good" -> bold, and the stray "_GoBack" bookmark
#    moves to sit between "Camera" and " should be good"
# ---------------------------------------------------------------------
$camPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Camera should be good*") {
        $camPara = $p
        break
    }
}
$camPara.Range.Font.Bold = $true
$splitPos = $camPara.Range.Start + [int]"Camera".Length
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))

# ---------------------------------------------------------------------
# 3) "Writing to file" -> strike-through (paragraph mark + run)
# ---------------------------------------------------------------------
$writingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Writing to file*") {
        $writingPara = $p
        break
    }
}
$writingPara.Range.Font.StrikeThrough = $true

# ---------------------------------------------------------------------
# 4) New bold bullet "Eligibility (make buttons / text easy to read)"
#    right after "Writing to file"
# ---------------------------------------------------------------------
$insertionPoint = $d.Range($writingPara.Range.End, $writingPara.Range.End)
$insertionPoint.InsertAfter("Eligibility (make buttons / text easy to read)" + [char]13)
$eligibilityPara = $writingPara.Next()
$eligibilityPara.Range.Font.Bold = $true

# ---------------------------------------------------------------------
# 5) "Changing speedometer to image of squid" -> bold
# ---------------------------------------------------------------------
$squidPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Changing speedometer to image of squid*") {
        $squidPara = $p
        break
    }
}
$squidPara.Range.Font.Bold = $true

# ---------------------------------------------------------------------
# 6) New bold bullet "Add switch code" right after "Changing speedometer..."
# ---------------------------------------------------------------------
$insertionPoint2 = $d.Range($squidPara.Range.End, $squidPara.Range.End)
$insertionPoint2.InsertAfter("Add switch code" + [char]13)
$switchPara = $squidPara.Next()
$switchPara.Range.Font.Bold = $true

# "Fixing toggle buttons" and "Do we really we want voltage to affect
# resistor thickness?" are left untouched (the latter simply loses the
# "_GoBack" bookmark, which already relocated in step 2).
